$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10 (TPM-recomputed values; the "ECs" sending cluster
# was dropped entirely, and Sending cluster values shift: FAPs/MuSCs/
# Resolving-Mac now occupy rows 2-10 in that order, 3 target rows each).
$rows = @(
    @("FAPs", "Mmp9", "Ephb2", "FAPs", 3, 1, 1.627176, 4.881528, 0.2586999172896354, 0.2586999172896354, 3, 1, 8.882730333333333, 26.648191, 0.9347132976570145, 0.9347132976570145, 14.453765612872, 130.083890515848, 0.241810252793392, 0.241810252793392),
    @("FAPs", "Mmp9", "Ephb2", "MuSCs", 3, 1, 1.627176, 4.881528, 0.2586999172896354, 0.2586999172896354, 3, 1, 0.616144, 1.848432, 0.06483569448352988, 0.0648356944835299, 1.002574729344, 9.023172564096001, 0.01677298880030525, 0.01677298880030525),
    @("FAPs", "Mmp9", "Ephb2", "Resolving-Mac", 3, 1, 1.627176, 4.881528, 0.2586999172896354, 0.2586999172896354, 2, 0.6666666666666666, 0.004286, 0.012858, 0.0004510078594555965, 0.0004510078594555965, 0.006974076336, 0.062766687024, 0.0001166756959381383, 0.0001166756959381383),
    @("MuSCs", "Mmp9", "Ephb2", "FAPs", 3, 1, 0.5265696666666667, 1.579709, 0.08371775961168156, 0.08371775961168154, 3, 1, 8.882730333333333, 26.648191, 0.9347132976570145, 0.9347132976570145, 4.677376350713223, 42.096387156419, 0.07825210315909209, 0.07825210315909208),
    @("MuSCs", "Mmp9", "Ephb2", "MuSCs", 3, 1, 0.5265696666666667, 1.579709, 0.08371775961168156, 0.08371775961168154, 3, 1, 0.616144, 1.848432, 0.06483569448352988, 0.0648356944835299, 0.3244427406986667, 2.919984666288, 0.005427899085028582, 0.005427899085028583),
    @("MuSCs", "Mmp9", "Ephb2", "Resolving-Mac", 3, 1, 0.5265696666666667, 1.579709, 0.08371775961168156, 0.08371775961168154, 2, 0.6666666666666666, 0.004286, 0.012858, 0.0004510078594555965, 0.0004510078594555965, 0.002256877591333333, 0.020311898322, 0.00003775736756088269, 0.00003775736756088269),
    @("Resolving-Mac", "Mmp9", "Ephb2", "FAPs", 3, 1, 4.136074666666667, 12.408224, 0.6575823230986831, 0.6575823230986831, 3, 1, 8.882730333333333, 26.648191, 0.9347132976570145, 0.9347132976570145, 36.73963590253155, 330.656723122784, 0.6146509417045305, 0.6146509417045305),
    @("Resolving-Mac", "Mmp9", "Ephb2", "MuSCs", 3, 1, 4.136074666666667, 12.408224, 0.6575823230986831, 0.6575823230986831, 3, 1, 0.616144, 1.848432, 0.06483569448352988, 0.0648356944835299, 2.548417589418667, 22.935758304768, 0.04263480659819605, 0.04263480659819607),
    @("Resolving-Mac", "Mmp9", "Ephb2", "Resolving-Mac", 3, 1, 4.136074666666667, 12.408224, 0.6575823230986831, 0.6575823230986831, 2, 0.6666666666666666, 0.004286, 0.012858, 0.0004510078594555965, 0.0004510078594555965, 0.01772721602133333, 0.159544944192, 0.0002965747959565755, 0.0002965747959565756)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

# Drop the old rows 11-13 (used to be the "Resolving-Mac" sender block,
# now folded into rows 8-10 above) so the table shrinks to A1:T10.
$ws.Range("A11:A13").EntireRow.Delete()
